# The deck currently ships two theme parts:
#   ppt/theme/theme1.xml -> generic "Office Theme" (used by the notes master)
#   ppt/theme/theme2.xml -> the custom "Integral" theme actually applied to the
#                            slide master / every slide
#
# The target revision swaps the two themes' colour palettes: the theme that is
# wired up as the presentation's live/active theme (theme2.xml, exposed by the
# object model as $p.SlideMaster.Theme / $p.NotesMaster.Theme / the per-slide
# ThemeColorScheme - they all resolve to the same live theme) should end up
# holding the plain "Office Theme" palette instead of "Integral".
#
# Recolour it one swatch at a time via ThemeColorScheme.Colors(i).RGB, which
# is the supported way to edit a live theme's palette through the PowerPoint
# object model.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeColor($scheme, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $ole = [System.Drawing.ColorTranslator]::ToOle([System.Drawing.Color]::FromArgb($r, $g, $b))
    $scheme.Colors($index).RGB = $ole
}

# Index order matches the a:clrScheme child order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
Set-ThemeColor $tcs 1  "000000"
Set-ThemeColor $tcs 2  "FFFFFF"
Set-ThemeColor $tcs 3  "44546A"
Set-ThemeColor $tcs 4  "E7E6E6"
Set-ThemeColor $tcs 5  "5B9BD5"
Set-ThemeColor $tcs 6  "ED7D31"
Set-ThemeColor $tcs 7  "A5A5A5"
Set-ThemeColor $tcs 8  "FFC000"
Set-ThemeColor $tcs 9  "4472C4"
Set-ThemeColor $tcs 10 "70AD47"
Set-ThemeColor $tcs 11 "0563C1"
Set-ThemeColor $tcs 12 "954F72"
